# Challenge_1_50days_tracker.xlsx
# Commit: "101. Symmetric Tree (DFS -recursion)"
#
# A new LeetCode entry (Q#101 "Symmetric Tree") is inserted into the
# tracker sheet right after row 77 ("226. Invert Binary Tree"), pushing
# every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78; this shifts rows 78..127 down to 79..128
# and inherits formatting/styles from the surrounding rows automatically.
$ws.Rows.Item(78).Insert()

# Populate the new row with the new problem's data.
$ws.Range("A78").Value2 = 101
$ws.Range("B78").Value2 = "Symmetric Tree"
$ws.Range("C78").Value2 = "Easy"
$ws.Range("D78").Value2 = "DFS ,recursion"
$ws.Range("E78").Value2 = 45839
$ws.Range("F78").Value2 = "Python"

# The existing "127. Word Ladder" hyperlink was anchored on the old B117;
# since the inserted row pushed that cell's content down to B118, re-anchor
# the hyperlink to follow it.
$ws.Range("B117").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B118"), "https://leetcode.com/problems/word-ladder/", "", "https://leetcode.com/problems/word-ladder/", "127. Word Ladder")

# Hyperlinks.Add overwrites the cell with Excel's default "Hyperlink" look;
# restore the sheet's original custom link styling (non-underlined, size 10,
# custom blue) that B118 ("127. Word Ladder") used before the edit.
$ws.Range("B118").Font.Size = 10
$ws.Range("B118").Font.Underline = 0
$ws.Range("B118").Font.Color = 16745482

# Reflect the author's final selection position in the saved workbook.
$ws.Range("F79").Select()
